$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 25: new data row (day 24 post)
$ws.Range("A25").Value = 24
$ws.Range("B25").NumberFormat = "d-mmm-yy"
$ws.Range("B25").Value = (Get-Date -Year 2020 -Month 1 -Day 24).Date
$ws.Range("C25").Value = "Surah Baqarah, 250 - 257"
$ws.Range("E25").Value = "Qasim Ali"
$ws.Range("F25").Value = "Tone of Quran, Philanthropist, Simple Islam, Book reading, Success through Quran"
$d25Text = @'
h1: What is wrong with us..?
p: I have read some books on personal development.
h6: 1. Awaken the giant within by <i>Anthony Robbins</i>
p: It changed my life back in 2016. I managed to quit smoking, get over with my divorce and become independant to take good care of my family.
h6: 2. Think and grow rich by <i>Napoleon Hill</i>
p: This book made it easier for me to pursue my goals. 
h6: 3. The 80 / 20 principle by <i>Richard Koch</i>
p: It made my days more productive. It gave me clarity of my days, told me to focus on the 20 % of the most productive part of the day.
h6: 4. Lean startup by <i>Eric Ries</i>
p: It made me iterate the cycle of validating, developing and launching products a little quicker. These blogs are a result of principle mentioned in this book. We have to iterate things faster.
h3: So… What is common in above books and Quran?
p: If somebody asks me what people who write such books on personal growth emphasise upon?
p: They all focus on some key factors:-
p.b-left: 1. One needs to have patience to achieve success.
p.b-left: 2. If one fails and does not get back on his feet, he can not succeed.
p.b-left: 3. Fight for victory, till one achieves it.
quote: And when they went forth to [face] Jaloot and his soldiers, they said, "Our Lord, pour upon us <span class=‘lavender’>patience</span> and <span class=‘lavender’>plant firmly our feet</span> and <span class=‘lavender’>give us victory</span> over the disbelieving people." <br> Surah Baqarah verse 250
p: Quran is also focusing on personal development. 
p: Why is it that we believe in all famous philanthropists who have written best sellers in past 1 year, but we can not believe in Prophet Muhammad (saw)?
p: Why can not we surrender to Allah’s will that He wants us to live a happy and contented life?
p: I think, the answer lies in our strengthened bond to self-made beliefs we have associated with Quran. Recently we were living with Hindus, before that corrupted muslim leaders and prior to that Mongols. Generation after generation we became a mix of traditions linked to these societies. We ended up detoriorating the concept of Quran.
p.b-left: Reading Quran needs ablution. Reading any other book is fine in bed, on floor, in car. But to read Quran we need to purify ourselves.
p.b-left: Learning Quran is recommended after establishing our footholds in society. Getting a good qualification, followed by a good job and finally getting married are considered pre-requisites to give time to Quran. Quran comes at 4th or 5th in priorities.
p.b-left: Acting upon Quran is labelled <b>“religious”</b>. Striving for teachings of Quran is labelled <b>”extremist”</b>. 
h3: Why cannot we consider Quran one Normal book to act upon?
p.note: Seriously guys, writing about Quran right now is so interesting. It is a bright book to learn.
p: Ok. Quran is a bright book. It needs to be understood in small parts, portions, ayats and words. Each time we read one verse from any where in the book, it talks different. Almost all verses talk about the might of our Protector. No one is a scholar of this Book. Some have read it more and some less. But each eye reading it gets its own meaning. You are going for marriage? It talks to relieve you of your stress. You are going for an interview? It talks of relying upon Allah. You are dying? It talks of life ahead. You are happy? It talks of becoming more happier. 
p: Just opening the book, opens so many avenues to explore. Reading it, hits our heart hard. Right after sinning, during repentance the hit is hardest. Later it fades and we end up sinning again. Reading books is a good habit. Reading Quran is a good habit just like all other good habits. Keeping the influx of encouraging sentences alive keeps us alive. Else what use is living, if it is just a piece of meat waking up and going back to sleep.
p: Quran is a light book. It teaches all those things, almost all succesful people talk about. If we can stick to “Think and grow rich” by <i>Napolean Hill</i>, it is equally important we consider “Quran” by <i>Allah (swt)</i> as good book to read before going to bed. One day read “Reader’s digest”, next day read “Quran” and third day read News paper. Make it part of routine affairs. Do not alienate this book.. It is just another book to keep us live a good life.
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
'@
$ws.Range("D25").Value = $d25Text

# Match row height of preceding long-content rows (max Excel row height)
$ws.Rows.Item(25).RowHeight = 409.6

# Move selection / scroll position to the new row, mirroring the author's saved view
$ws.Range("D25").Select() | Out-Null

Write-Host "Row 25 populated."
